$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted above the current row 154, pushing all
# subsequent rows (154-246) down by one (to 155-247). Insert a fresh row so
# the rest of the sheet shifts down automatically.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new record's data.
$ws.Range("A154").Value = 4
$ws.Range("B154").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C154").Value = "Los Lagos"
$ws.Range("D154").Value2 = 44704
$ws.Range("E154").Value = 10
$ws.Range("F154").Value = "Fruta"
$ws.Range("G154").Value = 100104
$ws.Range("H154").Value = "Frutos de pepita"
$ws.Range("I154").Value = 100104005
$ws.Range("J154").Value = "Pera"
$ws.Range("K154").Value = "Packham's Triumph"
$ws.Range("L154").Value = "Segunda"
$ws.Range("M154").Value = 300
$ws.Range("N154").Value = 12000
$ws.Range("O154").Value = 12000
$ws.Range("P154").Value = 12000
$ws.Range("Q154").Value = '$/caja 15 kilos empedrada'
$ws.Range("R154").Value = "Región de O'Higgins"
$ws.Range("S154").Value = 800
$ws.Range("T154").Value = 15
